$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Subject cell (B2) with the new value
$ws.Range("B2").Value = "Immuno&hema"

# Apply formatting to B2: centered alignment, light-gray fill, 11pt Calibri font
$r = $ws.Range("B2")
$r.HorizontalAlignment = -4108   # xlCenter
$r.VerticalAlignment = -4108     # xlCenter
$r.Interior.Color = 15790320     # RGB(240,240,240) = #F0F0F0
$r.Interior.PatternColor = 15790320
$r.Font.Size = 11

# Mirror the author's final selection landing on the edited cell
[void]$r.Select()
